# Scheduled runner refresh: pushes newly-pulled Market Board prices/profits
# (columns H:N) into each job sheet's leve table. A few rows also gain or
# lose a LeveProfit cell where the upstream feed started/stopped reporting it.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 833509.0600000001
$ws.Range("I6").Value = 1000030.9
$ws.Range("J6").Value = 900
$ws.Range("K6").Value = 3000092.7
$ws.Range("L6").Value = 2700
$ws.Range("M6").Value = -2999980.7
$ws.Range("N6").Value = -2924
$ws.Range("H76").Value = 4701
$ws.Range("I76").Value = 4103
$ws.Range("K76").Value = 4103
$ws.Range("M76").Value = -3788
$ws.Range("H79").Value = 4701
$ws.Range("I79").Value = 4103
$ws.Range("K79").Value = 4103
$ws.Range("M79").Value = -3011
$ws.Range("H112").Value = 984.16
$ws.Range("J112").Value = 1013.2174
$ws.Range("L112").Value = 3039.6522
$ws.Range("N112").Value = -5255.6522
$ws.Range("H137").Value = 1822.92
$ws.Range("I137").Value = 1503.5238
$ws.Range("J137").Value = 3499.75
$ws.Range("K137").Value = 4510.5714
$ws.Range("L137").Value = 10499.25
$ws.Range("M137").Value = -1960.5714
$ws.Range("N137").Value = -15599.25
$ws.Range("H138").Value = 1404.8572
$ws.Range("I138").Value = 1234.3243
$ws.Range("J138").Value = 2666.8
$ws.Range("K138").Value = 3702.9729
$ws.Range("L138").Value = 8000.400000000001
$ws.Range("M138").Value = 1437.0271
$ws.Range("N138").Value = -18280.4

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18743.959
$ws.Range("I32").Value = 3897.5632
$ws.Range("J32").Value = 126380.336
$ws.Range("K32").Value = 3897.5632
$ws.Range("L32").Value = 126380.336
$ws.Range("M32").Value = -3610.5632
$ws.Range("N32").Value = -126954.336
$ws.Range("H45").Value = 127457.375
$ws.Range("I45").Value = 167949.83
$ws.Range("K45").Value = 167949.83
$ws.Range("M45").Value = -167572.83
$ws.Range("H122").Value = 1961.8077
$ws.Range("I122").Value = 1995.3478
$ws.Range("J122").Value = 1704.6666
$ws.Range("K122").Value = 5986.0434
$ws.Range("L122").Value = 5113.9998
$ws.Range("M122").Value = -3536.0434
$ws.Range("N122").Value = -10013.9998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 409
$ws.Range("I22").Value = 391.25
$ws.Range("J22").Value = 444.5
$ws.Range("K22").Value = 391.25
$ws.Range("L22").Value = 444.5
$ws.Range("M22").Value = -218.25
$ws.Range("N22").Value = -790.5
$ws.Range("H86").Value = 40778.83
$ws.Range("I86").Value = 57957.8
$ws.Range("J86").Value = 2603.3333
$ws.Range("K86").Value = 57957.8
$ws.Range("L86").Value = 2603.3333
$ws.Range("M86").Value = -56834.8
$ws.Range("N86").Value = -4849.3333
$ws.Range("H89").Value = 40778.83
$ws.Range("I89").Value = 57957.8
$ws.Range("J89").Value = 2603.3333
$ws.Range("K89").Value = 289789
$ws.Range("L89").Value = 13016.6665
$ws.Range("M89").Value = -284173
$ws.Range("N89").Value = -24248.6665
$ws.Range("H96").Value = 28000
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H134").Value = 3721.5745
$ws.Range("I134").Value = 3408.1082
$ws.Range("J134").Value = 4881.4
$ws.Range("K134").Value = 10224.3246
$ws.Range("L134").Value = 14644.2
$ws.Range("M134").Value = -7689.3246
$ws.Range("N134").Value = -19714.2

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1197.5217
$ws.Range("I58").Value = 994.0952
$ws.Range("J58").Value = 3333.5
$ws.Range("K58").Value = 994.0952
$ws.Range("L58").Value = 3333.5
$ws.Range("M58").Value = -791.0952
$ws.Range("N58").Value = -3739.5
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H122").Value = 798.5
$ws.Range("J122").Value = 798.5
$ws.Range("L122").Value = 2395.5
$ws.Range("N122").Value = -7295.5
$ws.Range("H136").Value = 1197.5217
$ws.Range("I136").Value = 994.0952
$ws.Range("J136").Value = 3333.5
$ws.Range("K136").Value = 2982.2856
$ws.Range("L136").Value = 10000.5
$ws.Range("M136").Value = -432.2856000000002
$ws.Range("N136").Value = -15100.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 12272.444
$ws.Range("J80").Value = 12272.444
$ws.Range("L80").Value = 36817.33199999999
$ws.Range("N80").Value = -38689.33199999999
$ws.Range("H83").Value = 12272.444
$ws.Range("J83").Value = 12272.444
$ws.Range("L83").Value = 110451.996
$ws.Range("N83").Value = -119811.996
$ws.Range("H122").Value = 350.5
$ws.Range("J122").Value = 300
$ws.Range("L122").Value = 2700
$ws.Range("N122").Value = -7600

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3198.1765
$ws.Range("I102").Value = 3284.4285
$ws.Range("J102").Value = 3137.8
$ws.Range("K102").Value = 3284.4285
$ws.Range("L102").Value = 3137.8
$ws.Range("M102").Value = -1662.4285
$ws.Range("N102").Value = -6381.8
$ws.Range("H112").Value = 40799
$ws.Range("J112").Value = 40799
$ws.Range("L112").Value = 40799
$ws.Range("N112").Value = -43015

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1045.7428
$ws.Range("I22").Value = 674.9167
$ws.Range("J22").Value = 1239.2174
$ws.Range("K22").Value = 674.9167
$ws.Range("L22").Value = 1239.2174
$ws.Range("M22").Value = -379.9167
$ws.Range("N22").Value = -1829.2174
$ws.Range("H27").Value = 1045.7428
$ws.Range("I27").Value = 674.9167
$ws.Range("J27").Value = 1239.2174
$ws.Range("K27").Value = 674.9167
$ws.Range("L27").Value = 1239.2174
$ws.Range("M27").Value = -567.9167
$ws.Range("N27").Value = -1453.2174
$ws.Range("H40").Value = 54678.844
$ws.Range("I40").Value = 334166.66
$ws.Range("K40").Value = 334166.66
$ws.Range("M40").Value = -334030.66
$ws.Range("H55").Value = 392.65518
$ws.Range("J55").Value = 516
$ws.Range("L55").Value = 516
$ws.Range("N55").Value = -862
$ws.Range("H68").Value = 2899.6191
$ws.Range("I68").Value = 1775
$ws.Range("J68").Value = 3164.2354
$ws.Range("K68").Value = 1775
$ws.Range("L68").Value = 3164.2354
$ws.Range("M68").Value = -1026
$ws.Range("N68").Value = -4662.2354
$ws.Range("H71").Value = 2899.6191
$ws.Range("I71").Value = 1775
$ws.Range("J71").Value = 3164.2354
$ws.Range("K71").Value = 8875
$ws.Range("L71").Value = 15821.177
$ws.Range("M71").Value = -5131
$ws.Range("N71").Value = -23309.177
$ws.Range("H132").Value = 3581.7932
$ws.Range("I132").Value = 3599.375
$ws.Range("J132").Value = 3497.4
$ws.Range("K132").Value = 10798.125
$ws.Range("L132").Value = 10492.2
$ws.Range("M132").Value = -8268.125
$ws.Range("N132").Value = -15552.2
$ws.Range("H136").Value = 1684.1578
$ws.Range("I136").Value = 1428.375
$ws.Range("J136").Value = 3048.3333
$ws.Range("K136").Value = 4285.125
$ws.Range("L136").Value = 9144.999899999999
$ws.Range("M136").Value = -1735.125
$ws.Range("N136").Value = -14244.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2502789.2
$ws.Range("I62").Value = 6252798
$ws.Range("J62").Value = 2783.3333
$ws.Range("K62").Value = 6252798
$ws.Range("L62").Value = 2783.3333
$ws.Range("M62").Value = -6252174
$ws.Range("N62").Value = -4031.3333
$ws.Range("H65").Value = 2502789.2
$ws.Range("I65").Value = 6252798
$ws.Range("J65").Value = 2783.3333
$ws.Range("K65").Value = 31263990
$ws.Range("L65").Value = 13916.6665
$ws.Range("M65").Value = -31260870
$ws.Range("N65").Value = -20156.6665
$ws.Range("H132").Value = 2189.3333
$ws.Range("I132").Value = 2167.9575
$ws.Range("J132").Value = 2289.8
$ws.Range("K132").Value = 6503.872499999999
$ws.Range("L132").Value = 6869.400000000001
$ws.Range("M132").Value = -3973.872499999999
$ws.Range("N132").Value = -11929.4
$ws.Range("H136").Value = 732.8444
$ws.Range("I136").Value = 419.22858
$ws.Range("J136").Value = 1830.5
$ws.Range("K136").Value = 1257.68574
$ws.Range("L136").Value = 5491.5
$ws.Range("M136").Value = 1292.31426
$ws.Range("N136").Value = -10591.5

